# biofire_pdf_data.xlsx — populate real Covid/Flu A/Flu B/RSV result
# expressions for the auto PDF-filler (replacing the old placeholder
# "None" / raw-organism-list text in column E with concrete per-test
# result columns E:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 is the "Negative" sample -> all four assay columns read Negative.
$ws.Range("E2:H2").Value = "Negative"

# Row 3 is the "Positive" sample -> all four assay columns read Positive.
$ws.Range("E3:H3").Value = "Positive"

# Keep the workbook window geometry in sync with the editor that produced
# this revision (best-effort; window chrome, not worksheet content).
$excel.Left = -108
$excel.Top = -108
$excel.Width = 23256
$excel.Height = 12456
